$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06993827754423682
$ws.Range("C2").Value = 1.28051242846386
$ws.Range("D2").Value = 0.1118724862322458
$ws.Range("B3").Value = 0.07145923101128866
$ws.Range("C3").Value = 0.3653689161226725
$ws.Range("D3").Value = 0.08092646617548818
$ws.Range("B4").Value = 0.0254584143442825
$ws.Range("C4").Value = -0.4374107944482474
$ws.Range("D4").Value = 0.05605660395857231
$ws.Range("B5").Value = 0.05328850242883358
$ws.Range("C5").Value = -1.405751457105617
$ws.Range("D5").Value = 0.04873707937447077
$ws.Range("B6").Value = 0.05455910920632782
$ws.Range("C6").Value = -1.731039329317289
$ws.Range("D6").Value = 0.06947594556545653
$ws.Range("B7").Value = 0.04802847949661104
$ws.Range("C7").Value = -2.457473450745832
$ws.Range("D7").Value = 0.103691524082343
$ws.Range("B8").Value = 0.04959071786274524
$ws.Range("C8").Value = -2.424510454062424
$ws.Range("D8").Value = 0.08238748799846336
$ws.Range("B9").Value = 0.03230547907421665
$ws.Range("C9").Value = -2.159009737473113
$ws.Range("D9").Value = 0.05389533225787194
$ws.Range("B10").Value = 0.05015817493312774
$ws.Range("C10").Value = -1.684239037260034
$ws.Range("D10").Value = 0.1050997489056637
$ws.Range("B11").Value = 0.04429376798151233
$ws.Range("C11").Value = -0.9502772827632566
$ws.Range("D11").Value = 0.07611530644491771
$ws.Range("B12").Value = 0.05862732986404143
$ws.Range("C12").Value = -0.05127554389323999
$ws.Range("D12").Value = 0.07896955599256779
$ws.Range("B13").Value = 0.0590147319268373
$ws.Range("C13").Value = 0.8503318604047228
$ws.Range("D13").Value = 0.08671861254113698
$ws.Range("B14").Value = 0.05508826264160128
$ws.Range("C14").Value = 1.612384144361176
$ws.Range("D14").Value = 0.1000246706241896
$ws.Range("B15").Value = 0.0515852461788705
$ws.Range("C15").Value = 2.034225911906142
$ws.Range("D15").Value = 0.06439692709560274
$ws.Range("B16").Value = 0.05760902023758323
$ws.Range("C16").Value = 2.393022483076007
$ws.Range("D16").Value = 0.08356460243363614
$ws.Range("B17").Value = 0.03377036796550489
$ws.Range("C17").Value = 2.455401831136851
$ws.Range("D17").Value = 0.08362945325109637
$ws.Range("B18").Value = 0.06860553834659003
$ws.Range("C18").Value = 2.05303728248256
$ws.Range("D18").Value = 0.05120992943220692
$ws.Range("B19").Value = 0.02782889964574108
$ws.Range("C19").Value = 1.507401522273392
$ws.Range("D19").Value = 0.07083716993144563
$ws.Range("B20").Value = 0.0324205808074701
$ws.Range("C20").Value = 0.6965527366921982
$ws.Range("D20").Value = 0.07080628667109412
$ws.Range("B21").Value = 0.04727553853307719
$ws.Range("C21").Value = -0.1743236302332338
$ws.Range("D21").Value = 0.1003161406649846
$ws.Range("B22").Value = 0.0594780288679938
$ws.Range("C22").Value = -1.143159088146192
$ws.Range("D22").Value = 0.1002070098640018
$ws.Range("B23").Value = 0.05993713717158124
$ws.Range("C23").Value = -1.735473921028106
$ws.Range("D23").Value = 0.107584484309577
$ws.Range("B24").Value = 0.0438076884547291
$ws.Range("C24").Value = -2.299765500584122
$ws.Range("D24").Value = 0.08716557339378048
$ws.Range("B25").Value = 0.03055626529535757
$ws.Range("C25").Value = -2.655938911974428
$ws.Range("D25").Value = 0.07964231354411151
$ws.Range("B26").Value = 0.03839341754284745
$ws.Range("C26").Value = -2.337502576633018
$ws.Range("D26").Value = 0.0992286706519033
$ws.Range("B27").Value = 0.02560499860769954
$ws.Range("C27").Value = -2.059846620228276
$ws.Range("D27").Value = 0.07597648962284616
$ws.Range("B28").Value = 0.03164229339536411
$ws.Range("C28").Value = -1.140601647102034
$ws.Range("D28").Value = 0.0683117137390553
$ws.Range("B29").Value = 0.04861168616380662
$ws.Range("C29").Value = -0.4290268899286393
$ws.Range("D29").Value = 0.09224680711727345
$ws.Range("B30").Value = 0.03919435478712565
$ws.Range("C30").Value = 0.6280264612853667
$ws.Range("D30").Value = 0.06097825787215714
$ws.Range("B31").Value = 0.03973383891424129
$ws.Range("C31").Value = 1.413195814605684
$ws.Range("D31").Value = 0.06371145584359497
$ws.Range("B32").Value = 0.06200217437958402
$ws.Range("C32").Value = 2.087911496232874
$ws.Range("D32").Value = 0.1007856168527
$ws.Range("B33").Value = 0.05732318083683485
$ws.Range("C33").Value = 2.258732545881426
$ws.Range("D33").Value = 0.1122469837881425
$ws.Range("B34").Value = 0.07357778822896557
$ws.Range("C34").Value = 2.506916827043825
$ws.Range("D34").Value = 0.07272840995529302
$ws.Range("B35").Value = 0.04379520590964318
$ws.Range("C35").Value = 2.340303521220163
$ws.Range("D35").Value = 0.09632992271537255
$ws.Range("B36").Value = 0.03725080425939261
$ws.Range("C36").Value = 1.685188902104414
$ws.Range("D36").Value = 0.06217827690751524
$ws.Range("B37").Value = 0.05670827136956787
$ws.Range("C37").Value = 1.049108466779221
$ws.Range("D37").Value = 0.07982775483311866
$ws.Range("B38").Value = 0.02608827700342055
$ws.Range("C38").Value = -0.01160539147173432
$ws.Range("D38").Value = 0.07282107666081328
$ws.Range("B39").Value = 0.0610757242406389
$ws.Range("C39").Value = -0.7909789796802305
$ws.Range("D39").Value = 0.0490497442584437
$ws.Range("B40").Value = 0.03301426777284067
$ws.Range("C40").Value = -1.587783764468753
$ws.Range("D40").Value = 0.09175660600691256
$ws.Range("B41").Value = 0.03186259721330165
$ws.Range("C41").Value = -2.217847445097491
$ws.Range("D41").Value = 0.06319999942544281
